$d = $word.ActiveDocument

# Common namespace declarations used by the fragments we splice in via
# Range.InsertXML (InsertXML parses each fragment independently, so every
# namespace prefix used inside a fragment must be (re)declared on it).
$W_NS   = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$W14_NS = "xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'"
$DRAW_NS = "xmlns:wp='http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing' xmlns:wp14='http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing' xmlns:a='http://schemas.openxmlformats.org/drawingml/2006/main' xmlns:pic='http://schemas.openxmlformats.org/drawingml/2006/picture' xmlns:r='http://schemas.openxmlformats.org/officeDocument/2006/relationships'"

# ---------------------------------------------------------------------
# 1) Picture paragraph: mark the run holding the inline drawing as
#    <w:noProof/> (Word stamps this on image runs so the proofer skips
#    them).
# ---------------------------------------------------------------------
$picPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.InlineShapes.Count -gt 0) {
        $picPara = $i
        break
    }
}

$drawing = '<w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="3A3D3A18" wp14:editId="59EA6C81"><wp:extent cx="3667637" cy="1790950"/><wp:effectExtent l="0" t="0" r="9525" b="0"/><wp:docPr id="1349850739" name="Picture 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr id="1349850739" name=""/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId6"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="3667637" cy="1790950"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing>'
$picXml = "<w:p $W_NS $DRAW_NS><w:pPr><w:spacing w:line=`"240`" w:lineRule=`"auto`"/><w:jc w:val=`"center`"/></w:pPr><w:r w:rsidRPr=`"00062BC0`"><w:rPr><w:noProof/></w:rPr>$drawing</w:r></w:p>"
$d.Paragraphs($picPara).Range.InsertXML($picXml)

# ---------------------------------------------------------------------
# 2) "Hyperparameters" paragraph: split "Adam with lr 0.00005" and
#    "Early stopping with 2 epochs as patience" around proofing-error
#    markers, as a fresh round of spell/grammar check would.
# ---------------------------------------------------------------------
$hpPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Adam with lr*") {
        $hpPara = $i
        break
    }
}

$hpXml = @"
<w:p $W_NS>
  <w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>
  <w:r w:rsidRPr="003E08CF"><w:rPr><w:i/><w:iCs/></w:rPr><w:sym w:font="Wingdings" w:char="F0E0"/></w:r>
  <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Hyperparameters</w:t></w:r>
  <w:r w:rsidR="00891164"><w:rPr><w:i/><w:iCs/></w:rPr><w:br/></w:r>
  <w:r><w:t xml:space="preserve">Adam with </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>lr</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> 0.00005</w:t></w:r>
  <w:r w:rsidR="005725F6"><w:br/><w:t>Validation split of 0.1</w:t></w:r>
  <w:r w:rsidR="00DD6E1E"><w:br/><w:t>epochs: 50, batch size:16</w:t></w:r>
  <w:r w:rsidR="001A519F"><w:br/><w:t>Dropout in dense layer with a rate of 0.1</w:t></w:r>
  <w:r><w:br/><w:t xml:space="preserve">Early stopping with 2 epochs as </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>patience</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
</w:p>
"@
$d.Paragraphs($hpPara).Range.InsertXML($hpXml)

# ---------------------------------------------------------------------
# 3) "nlpaug, textblob, augly" paragraph: split into one run per tool
#    name wrapped in spellStart/spellEnd proofErr markers.
# ---------------------------------------------------------------------
$augPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*nlpaug*") {
        $augPara = $i
        break
    }
}

$rPr = '<w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Open Sans"/><w:color w:val="000000" w:themeColor="text1"/><w:kern w:val="0"/><w:sz w:val="26"/><w:szCs w:val="26"/><w14:ligatures w14:val="none"/></w:rPr>'
$augXml = @"
<w:p $W_NS $W14_NS>
  <w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="720"/><w:textAlignment w:val="baseline"/>$rPr</w:pPr>
  <w:r w:rsidRPr="00A91C10">$rPr<w:sym w:font="Wingdings" w:char="F0E0"/></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r w:rsidRPr="00A91C10">$rPr<w:t>nlpaug</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r w:rsidRPr="00A91C10">$rPr<w:t xml:space="preserve">, </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r w:rsidRPr="00A91C10">$rPr<w:t>textblob</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r w:rsidRPr="00A91C10">$rPr<w:t xml:space="preserve">, </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r w:rsidR="009E1E18" w:rsidRPr="00A91C10">$rPr<w:t>augly</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
"@
$d.Paragraphs($augPara).Range.InsertXML($augXml)

# ---------------------------------------------------------------------
# 4) Drop the trailing empty paragraph right before the section break,
#    after the final "Accuracy, precision, recall, F1 score" paragraph.
# ---------------------------------------------------------------------
$lastTextPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Accuracy, precision, recall, F1 score*") {
        $lastTextPara = $i
    }
}
$trailingPara = $lastTextPara + 1
$spanStart = $d.Paragraphs($lastTextPara).Range.Start
$spanEnd = $d.Paragraphs($trailingPara).Range.End
$spanRng = $d.Range($spanStart, $spanEnd)

$tailXml = "<w:p $W_NS><w:pPr><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t>Accuracy, precision, recall, F1 score</w:t></w:r></w:p>"
$spanRng.InsertXML($tailXml)

Write-Output "done; paragraphs=$($d.Paragraphs.Count)"
